$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 8).Value = 76.666664
$ws.Cells.Item(5, 9).Value = 76.666664
$ws.Cells.Item(5, 11).Value = 76.666664
$ws.Cells.Item(5, 13).Value = 38.333336
$ws.Cells.Item(28, 8).Value = 750
$ws.Cells.Item(28, 9).Value = 750
$ws.Cells.Item(28, 11).Value = 750
$ws.Cells.Item(28, 13).Value = -265
$ws.Cells.Item(76, 8).Value = 5348.2856
$ws.Cells.Item(76, 9).Value = 5798.8
$ws.Cells.Item(76, 10).Value = 4222
$ws.Cells.Item(76, 11).Value = 5798.8
$ws.Cells.Item(76, 12).Value = 4222
$ws.Cells.Item(76, 13).Value = -5483.8
$ws.Cells.Item(76, 14).Value = -4852
$ws.Cells.Item(79, 8).Value = 5348.2856
$ws.Cells.Item(79, 9).Value = 5798.8
$ws.Cells.Item(79, 10).Value = 4222
$ws.Cells.Item(79, 11).Value = 5798.8
$ws.Cells.Item(79, 12).Value = 4222
$ws.Cells.Item(79, 13).Value = -4706.8
$ws.Cells.Item(79, 14).Value = -6406
$ws.Cells.Item(96, 8).Value = 1092.9642
$ws.Cells.Item(96, 9).Value = 906.1667
$ws.Cells.Item(96, 10).Value = 1429.2
$ws.Cells.Item(96, 11).Value = 2718.5001
$ws.Cells.Item(96, 12).Value = 4287.6
$ws.Cells.Item(96, 13).Value = -1345.5001
$ws.Cells.Item(96, 14).Value = -7033.6
$ws.Cells.Item(113, 8).Value = 10492.607
$ws.Cells.Item(113, 9).Value = 11199.75
$ws.Cells.Item(113, 10).Value = 6249.75
$ws.Cells.Item(113, 11).Value = 11199.75
$ws.Cells.Item(113, 12).Value = 6249.75
$ws.Cells.Item(113, 13).Value = -7945.75
$ws.Cells.Item(113, 14).Value = -12757.75
$ws.Cells.Item(132, 8).Value = 4041.1936
$ws.Cells.Item(132, 9).Value = 2066.1
$ws.Cells.Item(132, 11).Value = 6198.299999999999
$ws.Cells.Item(132, 13).Value = -3668.299999999999
$ws.Cells.Item(135, 8).Value = 2743.762
$ws.Cells.Item(135, 9).Value = 1446.1875
$ws.Cells.Item(135, 11).Value = 13015.6875
$ws.Cells.Item(135, 13).Value = -10480.6875
$ws.Cells.Item(141, 8).Value = 5338
$ws.Cells.Item(141, 9).Value = 5277.875
$ws.Cells.Item(141, 10).Value = 6300
$ws.Cells.Item(141, 11).Value = 15833.625
$ws.Cells.Item(141, 12).Value = 18900
$ws.Cells.Item(141, 13).Value = -10653.625
$ws.Cells.Item(141, 14).Value = -29260
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 6497
$ws.Cells.Item(2, 9).Value = 8406.467000000001
$ws.Cells.Item(2, 11).Value = 8406.467000000001
$ws.Cells.Item(2, 13).Value = -8293.467000000001
$ws.Cells.Item(45, 8).Value = 1061.56
$ws.Cells.Item(45, 9).Value = 1039
$ws.Cells.Item(45, 10).Value = 1090.2727
$ws.Cells.Item(45, 11).Value = 1039
$ws.Cells.Item(45, 12).Value = 1090.2727
$ws.Cells.Item(45, 13).Value = -662
$ws.Cells.Item(45, 14).Value = -1844.2727
$ws.Cells.Item(74, 8).Value = 52379.1
$ws.Cells.Item(74, 9).Value = 58208.887
$ws.Cells.Item(74, 10).Value = 1368.5
$ws.Cells.Item(74, 11).Value = 58208.887
$ws.Cells.Item(74, 12).Value = 1368.5
$ws.Cells.Item(74, 13).Value = -57334.887
$ws.Cells.Item(74, 14).Value = -3116.5
$ws.Cells.Item(77, 8).Value = 52379.1
$ws.Cells.Item(77, 9).Value = 58208.887
$ws.Cells.Item(77, 10).Value = 1368.5
$ws.Cells.Item(77, 11).Value = 291044.435
$ws.Cells.Item(77, 12).Value = 6842.5
$ws.Cells.Item(77, 13).Value = -286676.435
$ws.Cells.Item(77, 14).Value = -15578.5
$ws.Cells.Item(101, 8).Value = 19999
$ws.Cells.Item(101, 10).Value = 19999
$ws.Cells.Item(101, 12).Value = 19999
$ws.Cells.Item(101, 14).Value = -26489
$ws.Cells.Item(116, 8).Value = 6497
$ws.Cells.Item(116, 9).Value = 8406.467000000001
$ws.Cells.Item(116, 11).Value = 8406.467000000001
$ws.Cells.Item(116, 13).Value = -6112.467000000001
$ws.Cells.Item(132, 8).Value = 3573.342
$ws.Cells.Item(132, 9).Value = 3588.4
$ws.Cells.Item(132, 11).Value = 10765.2
$ws.Cells.Item(132, 13).Value = -8235.200000000001
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 6497
$ws.Cells.Item(3, 9).Value = 8406.467000000001
$ws.Cells.Item(3, 11).Value = 8406.467000000001
$ws.Cells.Item(3, 13).Value = -8292.467000000001
$ws.Cells.Item(134, 8).Value = 1765.3143
$ws.Cells.Item(134, 9).Value = 1272.4814
$ws.Cells.Item(134, 11).Value = 3817.4442
$ws.Cells.Item(134, 13).Value = -1282.4442
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 8).Value = 142861120
$ws.Cells.Item(4, 9).Value = 4000
$ws.Cells.Item(4, 10).Value = 200003980
$ws.Cells.Item(4, 11).Value = 4000
$ws.Cells.Item(4, 12).Value = 200003980
$ws.Cells.Item(4, 13).Value = -3888
$ws.Cells.Item(4, 14).Value = -200004204
$ws.Cells.Item(7, 8).Value = 166666750
$ws.Cells.Item(7, 9).Value = 84.5
$ws.Cells.Item(7, 11).Value = 84.5
$ws.Cells.Item(7, 13).Value = 28.5
$ws.Cells.Item(14, 8).Value = 2998.6667
$ws.Cells.Item(14, 10).Value = 3998.5
$ws.Cells.Item(14, 12).Value = 3998.5
$ws.Cells.Item(14, 14).Value = -4338.5
$ws.Cells.Item(31, 8).Value = 149526.92
$ws.Cells.Item(31, 9).Value = 194212.9
$ws.Cells.Item(31, 10).Value = 38875.906
$ws.Cells.Item(31, 11).Value = 194212.9
$ws.Cells.Item(31, 12).Value = 38875.906
$ws.Cells.Item(31, 13).Value = -193917.9
$ws.Cells.Item(31, 14).Value = -39465.906
$ws.Cells.Item(34, 8).Value = 149526.92
$ws.Cells.Item(34, 9).Value = 194212.9
$ws.Cells.Item(34, 10).Value = 38875.906
$ws.Cells.Item(34, 11).Value = 194212.9
$ws.Cells.Item(34, 12).Value = 38875.906
$ws.Cells.Item(34, 13).Value = -194010.9
$ws.Cells.Item(34, 14).Value = -39279.906
$ws.Cells.Item(56, 8).Value = 46989
$ws.Cells.Item(56, 10).Value = 46989
$ws.Cells.Item(56, 12).Value = 46989
$ws.Cells.Item(56, 14).Value = -48679
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 991.7857
$ws.Cells.Item(122, 9).Value = 1027.9166
$ws.Cells.Item(122, 11).Value = 3083.7498
$ws.Cells.Item(122, 13).Value = -633.7498000000001
$ws.Cells.Item(141, 8).Value = 112277.98
$ws.Cells.Item(141, 10).Value = 112277.98
$ws.Cells.Item(141, 12).Value = 112277.98
$ws.Cells.Item(141, 14).Value = -122637.98
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(75, 8).Value = 2858.7856
$ws.Cells.Item(75, 9).Value = 2398.5
$ws.Cells.Item(75, 10).Value = 3042.9
$ws.Cells.Item(75, 11).Value = 7195.5
$ws.Cells.Item(75, 12).Value = 9128.700000000001
$ws.Cells.Item(75, 13).Value = -6197.5
$ws.Cells.Item(75, 14).Value = -11124.7
$ws.Cells.Item(78, 8).Value = 2858.7856
$ws.Cells.Item(78, 9).Value = 2398.5
$ws.Cells.Item(78, 10).Value = 3042.9
$ws.Cells.Item(78, 11).Value = 21586.5
$ws.Cells.Item(78, 12).Value = 27386.1
$ws.Cells.Item(78, 13).Value = -16594.5
$ws.Cells.Item(78, 14).Value = -37370.10000000001
$ws.Cells.Item(94, 8).Value = 800
$ws.Cells.Item(94, 10).Value = 800
$ws.Cells.Item(94, 12).Value = 2400
$ws.Cells.Item(94, 14).Value = -3752
$ws.Cells.Item(131, 8).Value = 6945973.5
$ws.Cells.Item(131, 10).Value = 1593.0308
$ws.Cells.Item(131, 12).Value = 4779.0924
$ws.Cells.Item(131, 14).Value = -14859.0924
$ws.Cells.Item(136, 8).Value = 1529.091
$ws.Cells.Item(136, 9).Value = 1529.091
$ws.Cells.Item(136, 11).Value = 4587.272999999999
$ws.Cells.Item(136, 13).Value = 512.7270000000008
$ws.Cells.Item(140, 8).Value = 5320928.5
$ws.Cells.Item(140, 9).Value = 25001088
$ws.Cells.Item(140, 11).Value = 75003264
$ws.Cells.Item(140, 13).Value = -74998084
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 2935.8
$ws.Cells.Item(102, 9).Value = 2935.8
$ws.Cells.Item(102, 11).Value = 2935.8
$ws.Cells.Item(102, 13).Value = -1313.8
$ws.Cells.Item(132, 8).Value = 84499.64
$ws.Cells.Item(132, 9).Value = 102542.55
$ws.Cells.Item(132, 11).Value = 307627.65
$ws.Cells.Item(132, 13).Value = -305097.65
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(46, 8).Value = 2670.5
$ws.Cells.Item(46, 9).Value = 1849.6666
$ws.Cells.Item(46, 11).Value = 1849.6666
$ws.Cells.Item(46, 13).Value = -1661.6666
$ws.Cells.Item(122, 8).Value = 6125.1665
$ws.Cells.Item(122, 9).Value = 6088
$ws.Cells.Item(122, 10).Value = 6199.5
$ws.Cells.Item(122, 11).Value = 18264
$ws.Cells.Item(122, 12).Value = 18598.5
$ws.Cells.Item(122, 13).Value = -15814
$ws.Cells.Item(122, 14).Value = -23498.5
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(70, 8).Value = 37775.5
$ws.Cells.Item(70, 9).Value = 74999
$ws.Cells.Item(70, 10).Value = 25367.666
$ws.Cells.Item(70, 11).Value = 74999
$ws.Cells.Item(70, 12).Value = 25367.666
$ws.Cells.Item(70, 13).Value = -74684
$ws.Cells.Item(70, 14).Value = -25997.666
$ws.Cells.Item(73, 8).Value = 37775.5
$ws.Cells.Item(73, 9).Value = 74999
$ws.Cells.Item(73, 10).Value = 25367.666
$ws.Cells.Item(73, 11).Value = 74999
$ws.Cells.Item(73, 12).Value = 25367.666
$ws.Cells.Item(73, 13).Value = -73907
$ws.Cells.Item(73, 14).Value = -27551.666
$ws.Cells.Item(81, 8).Value = 3508.6667
$ws.Cells.Item(81, 9).Value = 3736.7273
$ws.Cells.Item(81, 11).Value = 7473.4546
$ws.Cells.Item(81, 13).Value = -6412.4546
$ws.Cells.Item(84, 8).Value = 3508.6667
$ws.Cells.Item(84, 9).Value = 3736.7273
$ws.Cells.Item(84, 11).Value = 37367.273
$ws.Cells.Item(84, 13).Value = -32063.273
$ws.Cells.Item(122, 8).Value = 3642.4285
$ws.Cells.Item(122, 9).Value = 3642.4285
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 10927.2855
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -8477.2855
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 84857
$ws.Cells.Item(135, 10).Value = 84857
$ws.Cells.Item(135, 12).Value = 84857
$ws.Cells.Item(135, 14).Value = -94997
$ws.Cells.Item(136, 8).Value = 215715.23
$ws.Cells.Item(136, 9).Value = 272736.53
$ws.Cells.Item(136, 11).Value = 818209.5900000001
$ws.Cells.Item(136, 13).Value = -815659.5900000001
